$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) so numeric-looking strings (e.g. "309.21")
# are written back as text, matching the source data which stores every
# price/volume cell as a string rather than a number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.897.70'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.812.40'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '309.21'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").Value = '0.3653'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").Value = '0.07333'
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").Value = '0.8666'
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("D11").Value = '20.28'
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '1.811.80'
$ws.Range("E12").Value = '  -4.76%  '
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '0.07078'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '6.501'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '91.58'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '0.000008691'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '14.62'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '26.913.74'
$ws.Range("D22").Value = '5.290'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '2.031.81'
$ws.Range("E24").Value = '  -4.72%  '
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '150.46'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").Value = '18.25'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").Value = '2.152'
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").Value = '5.266'
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("D30").Value = '115.46'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '0.08932'
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").Value = '0.7529'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D33").Value = '1.155'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = '4.482'
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '2.912'
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = '1.083'
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("D38").Value = '0.05278'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").Value = '2.986'
$ws.Range("E39").Value = '  +1.86%  '
$ws.Range("D40").Value = '0.01946'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '7.190'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '0.5289'
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").Value = '2.281'
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("D44").Value = '0.1649'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").Value = '8.383'
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  -2.70%  '
$ws.Range("D47").Value = '10.40'
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '102.98'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").Value = '1.656'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '0.06291'
$ws.Range("E51").Value = '  +0.12%  '

# Restore the default (Normal) style on column D so no stray number-format
# style gets left on the cells themselves.
$ws.Range("D2:D51").Style = "Normal"

